# Apply the "completed finishing with the correct calculations and t stat"
# edit to stroopdata.xlsx.
#
# Before:
#   I31 = AVERAGE(I2:I25)
#   A36 = "Standard Error"   B36 = I31/SQRT(B27)
#   A37 = (blank)            B37 = (blank)
#   A38 = "T Stat"            B38 = B35/B36
#   L46 = footnote
#
# After:
#   I31 = SUM(I2:I25)/(COUNT(I2:I25)-1)
#   A36 = "STD of Different"  B36 = SQRT(I31)
#   A37 = "Standard Error"    B37 = B36/SQRT(B27)
#   A38 = (blank)             B38 = (blank)
#   A39 = "T Stat"            B39 = B35/B37
#   L47 = footnote

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the I31 aggregate: AVERAGE -> SUM/(COUNT-1) -------------------
$ws.Range("I31").Formula = "=SUM(I2:I25)/(COUNT(I2:I25)-1)"

# --- 2. Insert a fresh row above the old row 37 so everything below (the
#        blank spacer row, the T Stat row and the footnote further down)
#        shifts down by one, picking up matching formatting on the way. ---
$ws.Rows.Item(37).Insert()

# --- 3. Row 36 becomes "STD of Different" = SQRT(I31) ---------------------
$ws.Range("A36").Value = "STD of Different"
$ws.Range("B36").Formula = "=SQRT(I31)"

# --- 4. The newly inserted row 37 becomes "Standard Error" ----------------
$ws.Range("A37").Value = "Standard Error"
$ws.Range("B37").Formula = "=B36/SQRT(B27)"

# --- 5. Row 39 (old row 38, shifted down) now divides by the new Standard
#        Error row instead of the renamed Standard Deviation row. ----------
$ws.Range("B39").Formula = "=B35/B37"

# --- 6. Match the author's final selection/scroll state -------------------
$ws.Range("A27:E39").Select()
